# Generate Report for Handoff
# Updates the localization-status report: flips the "Handed back: in sync
# with en-US" status to "Ready for handoff" (a new handoff round has
# started) and refreshes the related timestamps, then narrows the
# "Status"/language result columns that used to be sized for the long
# status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Refreshed timestamps
$overview.Range("G2").Value = "2016-08-12 13:11:15"
$dede.Range("H2").Value = "2016-08-12 13:11:15"
$zhcn.Range("H2").Value = "2016-08-12 13:11:06"

# --- Narrow the columns that previously fit the long status text
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
